$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.427.62"
$ws.Range("E2").Value = "  -3.10%  "

$ws.Range("D3").Value = "'3.317.26"
$ws.Range("E3").Value = "  -3.10%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'558.56"
$ws.Range("E5").Value = "  -3.13%  "

$ws.Range("D6").Value = "'143.34"
$ws.Range("E6").Value = "  -3.34%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "'3.320.17"
$ws.Range("E8").Value = "  -2.98%  "

$ws.Range("D9").Value = "'0.478"
$ws.Range("E9").Value = "  -1.64%  "

$ws.Range("D10").Value = "'7.84"
$ws.Range("E10").Value = "  -2.05%  "

$ws.Range("D11").Value = "'0.119"
$ws.Range("E11").Value = "  -3.11%  "

$ws.Range("D12").Value = "'0.410"
$ws.Range("E12").Value = "  -0.77%  "

$ws.Range("D13").Value = "'3.883.54"
$ws.Range("E13").Value = "  -3.65%  "

$ws.Range("D15").Value = "'27.10"
$ws.Range("E15").Value = "  -3.94%  "

$ws.Range("D16").Value = "'3.313.10"
$ws.Range("E16").Value = "  -3.36%  "

$ws.Range("D17").Value = "'0.0000166"
$ws.Range("E17").Value = "  -3.03%  "

$ws.Range("D18").Value = "'60.406.20"
$ws.Range("E18").Value = "  -3.48%  "

$ws.Range("D19").Value = "'6.18"
$ws.Range("E19").Value = "  -2.96%  "

$ws.Range("D20").Value = "'14.53"
$ws.Range("E20").Value = "  +0.81%  "

$ws.Range("D21").Value = "'8.64"
$ws.Range("E21").Value = "  -3.07%  "

$ws.Range("D22").Value = "'374.98"
$ws.Range("E22").Value = "  -2.11%  "

$ws.Range("D23").Value = "'74.15"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("D24").Value = "'0.545"
$ws.Range("E24").Value = "  -3.65%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").Value = "'3.436.40"
$ws.Range("E26").Value = "  -4.31%  "

$ws.Range("E27").Value = "  -7.06%  "

$ws.Range("E28").Value = "  -5.80%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").Value = "'7.29"
$ws.Range("E30").Value = "  -4.32%  "

$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("D32").Value = "'7.70"
$ws.Range("E32").Value = "  -2.88%  "

$ws.Range("E33").Value = "  -3.35%  "

$ws.Range("D34").Value = "'22.60"
$ws.Range("E34").Value = "  -2.22%  "

$ws.Range("E35").Value = "  -4.16%  "

$ws.Range("D36").Value = "'5.22"
$ws.Range("E36").Value = "  -3.66%  "

$ws.Range("D37").Value = "'1.55"
$ws.Range("E37").Value = "  -5.16%  "

$ws.Range("D38").Value = "'166.60"
$ws.Range("E38").Value = "  -1.34%  "

$ws.Range("E39").Value = "  -1.80%  "

$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").Value = "'3.346.31"
$ws.Range("E40").Value = "  -3.63%  "

$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'26.87"
$ws.Range("E41").Value = "  -15.59%  "

$ws.Range("D42").Value = "'0.0743"
$ws.Range("E42").Value = "  -4.87%  "

$ws.Range("D43").Value = "'42.04"
$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("D44").Value = "'0.754"
$ws.Range("E44").Value = "  -3.43%  "

$ws.Range("D45").Value = "'4.21"
$ws.Range("E45").Value = "  -3.49%  "

$ws.Range("D46").Value = "'1.60"
$ws.Range("E46").Value = "  -5.05%  "

$ws.Range("E47").Value = "  -3.46%  "

$ws.Range("D48").Value = "'2.365.22"
$ws.Range("E48").Value = "  -6.90%  "

$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  -0.22%  "

$ws.Range("D50").Value = "'6.56"
$ws.Range("E50").Value = "  -4.86%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'21.52"
$ws.Range("E51").Value = "  -4.11%  "
